$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.75
$ws.Range("C2").Value = 0.4615384615384616
$ws.Range("D2").Value = 0.5714285714285714
$ws.Range("E2").Value = 13
$ws.Range("B3").Value = 0.5625
$ws.Range("C3").Value = 0.8181818181818182
$ws.Range("D3").Value = 0.6666666666666666
$ws.Range("E3").Value = 11
$ws.Range("B4").Value = 0.625
$ws.Range("C4").Value = 0.625
$ws.Range("D4").Value = 0.625
$ws.Range("E4").Value = 0.625
$ws.Range("B5").Value = 0.65625
$ws.Range("C5").Value = 0.6398601398601399
$ws.Range("D5").Value = 0.6190476190476191
$ws.Range("E5").Value = 24
$ws.Range("B6").Value = 0.6640625
$ws.Range("C6").Value = 0.625
$ws.Range("D6").Value = 0.615079365079365
$ws.Range("E6").Value = 24
$ws.Range("B7").Value = 0.4705882352941176
$ws.Range("C7").Value = 0.6153846153846154
$ws.Range("D7").Value = 0.5333333333333333
$ws.Range("E7").Value = 13
$ws.Range("B8").Value = 0.2857142857142857
$ws.Range("C8").Value = 0.1818181818181818
$ws.Range("D8").Value = 0.2222222222222222
$ws.Range("E8").Value = 11
$ws.Range("B9").Value = 0.4166666666666667
$ws.Range("C9").Value = 0.4166666666666667
$ws.Range("D9").Value = 0.4166666666666667
$ws.Range("E9").Value = 0.4166666666666667
$ws.Range("B10").Value = 0.3781512605042017
$ws.Range("C10").Value = 0.3986013986013986
$ws.Range("D10").Value = 0.3777777777777778
$ws.Range("E10").Value = 24
$ws.Range("B11").Value = 0.3858543417366946
$ws.Range("C11").Value = 0.4166666666666667
$ws.Range("D11").Value = 0.3907407407407408
$ws.Range("E11").Value = 24
$ws.Range("B12").Value = 0.6
$ws.Range("C12").Value = 0.6923076923076923
$ws.Range("D12").Value = 0.6428571428571429
$ws.Range("E12").Value = 13
$ws.Range("B13").Value = 0.5555555555555556
$ws.Range("C13").Value = 0.4545454545454545
$ws.Range("D13").Value = 0.5
$ws.Range("E13").Value = 11
$ws.Range("B14").Value = 0.5833333333333334
$ws.Range("C14").Value = 0.5833333333333334
$ws.Range("D14").Value = 0.5833333333333334
$ws.Range("E14").Value = 0.5833333333333334
$ws.Range("B15").Value = 0.5777777777777777
$ws.Range("C15").Value = 0.5734265734265734
$ws.Range("D15").Value = 0.5714285714285714
$ws.Range("E15").Value = 24
$ws.Range("B16").Value = 0.5796296296296296
$ws.Range("C16").Value = 0.5833333333333334
$ws.Range("D16").Value = 0.5773809523809524
$ws.Range("E16").Value = 24
$ws.Range("B17").Value = 0.4615384615384616
$ws.Range("C17").Value = 0.4615384615384616
$ws.Range("D17").Value = 0.4615384615384616
$ws.Range("E17").Value = 13
$ws.Range("B18").Value = 0.3636363636363636
$ws.Range("C18").Value = 0.3636363636363636
$ws.Range("D18").Value = 0.3636363636363636
$ws.Range("E18").Value = 11
$ws.Range("B20").Value = 0.4125874125874126
$ws.Range("C20").Value = 0.4125874125874126
$ws.Range("D20").Value = 0.4125874125874126
$ws.Range("E20").Value = 24
$ws.Range("B21").Value = 0.4166666666666667
$ws.Range("C21").Value = 0.4166666666666667
$ws.Range("D21").Value = 0.4166666666666667
$ws.Range("E21").Value = 24
$ws.Range("B22").Value = 0.6111111111111112
$ws.Range("C22").Value = 0.8461538461538461
$ws.Range("D22").Value = 0.7096774193548387
$ws.Range("E22").Value = 13
$ws.Range("B23").Value = 0.6666666666666666
$ws.Range("C23").Value = 0.3636363636363636
$ws.Range("D23").Value = 0.4705882352941177
$ws.Range("E23").Value = 11
$ws.Range("B24").Value = 0.625
$ws.Range("C24").Value = 0.625
$ws.Range("D24").Value = 0.625
$ws.Range("E24").Value = 0.625
$ws.Range("B25").Value = 0.6388888888888888
$ws.Range("C25").Value = 0.6048951048951049
$ws.Range("D25").Value = 0.5901328273244782
$ws.Range("E25").Value = 24
$ws.Range("B26").Value = 0.6365740740740741
$ws.Range("C26").Value = 0.625
$ws.Range("D26").Value = 0.6000948766603416
$ws.Range("E26").Value = 24
